# Week 13 logging update for Target Depth Data (Bears)
# Updates the "H" (home) row stats on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 267
$wsOff.Range("C2").Value = 178
$wsOff.Range("D2").Value = 68
$wsOff.Range("E2").Value = 34
$wsOff.Range("F2").Value = 8

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 376
$wsDef.Range("C2").Value = 256
$wsDef.Range("D2").Value = 91
$wsDef.Range("E2").Value = 44
